# LeetCodeStats.xlsx update: add a new "Minimum Depth of Binary Tree" entry
# as row 26 on Sheet1 (day 24), with its hyperlink + url note, and move the
# active selection to E29 (Sheet2's summary row recalculates automatically
# since it's SUM/COUNT formulas over Sheet1 columns C:H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 26 data ------------------------------------------------------
$ws.Range("B26").Value = "Minimum Depth of Binary Tree"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 440
$ws.Range("F26").Value = 0.66
$ws.Range("G26").Value = 58.2
$ws.Range("H26").Value = 0.752
$ws.Range("I26").Value = "https://leetcode.com/problems/minimum-depth-of-binary-tree/submissions/1065754632/"

# --- Hyperlink on B26 ---------------------------------------------------
# Hyperlinks.Add sets the cell's displayed text to the hyperlink target as
# a side effect, so restore the problem-name text afterwards.
$h = $ws.Hyperlinks.Add($ws.Range("B26"), "https://leetcode.com/problems/minimum-depth-of-binary-tree")
$h.TextToDisplay() = "https://leetcode.com/problems/minimum-depth-of-binary-tree"
$ws.Range("B26").Value = "Minimum Depth of Binary Tree"

# Match the same cell style used by the other problem-name/hyperlink cells
# above it (B22:B25) instead of the ad-hoc style Hyperlinks.Add creates.
$ws.Range("B25").Copy()
$ws.Range("B26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- View state ---------------------------------------------------------
$ws.Range("E29").Select() | Out-Null
